# SponsorSynq Competitor Comparison — "Action Items" sheet update
#
# Removes the three completed High-Priority rows (Sponsor Matching,
# Verification System, Dashboard Savings Display) and the now-superseded
# "Ambassador Visibility" row, so the five remaining High-Priority items
# (Event Collaboration/Co-hosting, Promoter Referral System, Landing Page
# Overhaul, Venue Partnership System, Revenue Stream Documentation) slide
# up to directly follow the header row. Everything below (Medium/Low
# priority blocks) shifts up by the same 4 rows automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action Items")

# Rows 10-13 hold: ✓ Sponsor Matching, ✓ Verification System,
# ✓ Dashboard Savings Display, Ambassador Visibility.
# Deleting the whole rows shifts rows 14+ up by 4, preserving each row's
# own formatting/height/merges (matches the diff exactly: old row 14
# "Event Collaboration/Co-hosting" becomes the new row 10, etc.).
$ws.Range("A10:D13").EntireRow.Delete()
